# Auto-generated edit script: updates Leve profit calculation values
# across multiple worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 103
$ws.Range("H103").Value = 5312.4546
$ws.Range("I103").Value = 50000
$ws.Range("J103").Value = 843.7
$ws.Range("K103").Value = 150000
$ws.Range("L103").Value = 2531.1
$ws.Range("M103").Value = -149414
$ws.Range("N103").Value = -3703.1
# Row 106
$ws.Range("H106").Value = 484151.56
$ws.Range("I106").Value = 2005608.6
$ws.Range("J106").Value = 8696.25
$ws.Range("K106").Value = 2005608.6
$ws.Range("L106").Value = 8696.25
$ws.Range("M106").Value = -2004977.6
$ws.Range("N106").Value = -9958.25
# Row 113
$ws.Range("H113").Value = 4159.4595
$ws.Range("I113").Value = 3443.75
$ws.Range("J113").Value = 4704.7617
$ws.Range("K113").Value = 3443.75
$ws.Range("L113").Value = 4704.7617
$ws.Range("M113").Value = -189.75
$ws.Range("N113").Value = -11212.7617
# Row 137
$ws.Range("H137").Value = 2610.7273
$ws.Range("I137").Value = 3074
$ws.Range("J137").Value = 2290
$ws.Range("K137").Value = 9222
$ws.Range("L137").Value = 6870
$ws.Range("M137").Value = -6672
$ws.Range("N137").Value = -11970

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 23
$ws.Range("H23").Value = 23000
$ws.Range("J23").Value = 23000
$ws.Range("L23").Value = 23000
$ws.Range("N23").Value = -23518
# Row 32
$ws.Range("H32").Value = 8135.143
$ws.Range("I32").Value = 4162.1724
$ws.Range("J32").Value = 20263.158
$ws.Range("K32").Value = 4162.1724
$ws.Range("L32").Value = 20263.158
$ws.Range("M32").Value = -3875.1724
$ws.Range("N32").Value = -20837.158
# Row 45
$ws.Range("H45").Value = 9092205
$ws.Range("I45").Value = 10102228
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 10102228
$ws.Range("L45").Value = 2000
$ws.Range("M45").Value = -10101851
$ws.Range("N45").Value = -2754
# Row 122
$ws.Range("H122").Value = 1314.7742
$ws.Range("I122").Value = 1187.7858
$ws.Range("K122").Value = 3563.3574
$ws.Range("M122").Value = -1113.3574

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 4018.9143
$ws.Range("I134").Value = 4069.9355
$ws.Range("J134").Value = 3623.5
$ws.Range("K134").Value = 12209.8065
$ws.Range("L134").Value = 10870.5
$ws.Range("M134").Value = -9674.806500000001
$ws.Range("N134").Value = -15940.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2430.7046
$ws.Range("I31").Value = 1428.125
$ws.Range("J31").Value = 2653.5
$ws.Range("K31").Value = 1428.125
$ws.Range("L31").Value = 2653.5
$ws.Range("M31").Value = -1133.125
$ws.Range("N31").Value = -3243.5
# Row 34
$ws.Range("H34").Value = 2430.7046
$ws.Range("I34").Value = 1428.125
$ws.Range("J34").Value = 2653.5
$ws.Range("K34").Value = 1428.125
$ws.Range("L34").Value = 2653.5
$ws.Range("M34").Value = -1226.125
$ws.Range("N34").Value = -3057.5
# Row 99
$ws.Range("H99").Value = 2090.138
$ws.Range("I99").Value = 1880.9524
$ws.Range("J99").Value = 2639.25
$ws.Range("K99").Value = 1880.9524
$ws.Range("L99").Value = 2639.25
$ws.Range("M99").Value = -382.9523999999999
$ws.Range("N99").Value = -5635.25
# Row 105
$ws.Range("H105").Value = 861.63635
$ws.Range("I105").Value = 867.8
$ws.Range("J105").Value = 800
$ws.Range("K105").Value = 867.8
$ws.Range("L105").Value = 800
$ws.Range("M105").Value = 879.2
$ws.Range("N105").Value = -4294
# Row 122
$ws.Range("H122").Value = 1000
$ws.Range("I122").Value = 1000
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3000
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -550
# Row 126
$ws.Range("H126").Value = 2090.138
$ws.Range("I126").Value = 1880.9524
$ws.Range("J126").Value = 2639.25
$ws.Range("K126").Value = 5642.857199999999
$ws.Range("L126").Value = 7917.75
$ws.Range("M126").Value = -3172.857199999999
$ws.Range("N126").Value = -12857.75

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1193.5333
$ws.Range("I5").Value = 260.42856
$ws.Range("J5").Value = 2010
$ws.Range("K5").Value = 781.28568
$ws.Range("L5").Value = 6030
$ws.Range("M5").Value = -669.28568
$ws.Range("N5").Value = -6254
# Row 109
$ws.Range("H109").Value = 4989.364
$ws.Range("I109").Value = 1576
$ws.Range("J109").Value = 6269.375
$ws.Range("K109").Value = 4728
$ws.Range("L109").Value = 18808.125
$ws.Range("M109").Value = -3688
$ws.Range("N109").Value = -20888.125
# Row 132
$ws.Range("H132").Value = 2426.5925
$ws.Range("I132").Value = 2688
$ws.Range("J132").Value = 2099.8333
$ws.Range("K132").Value = 24192
$ws.Range("L132").Value = 18898.4997
$ws.Range("M132").Value = -21662
$ws.Range("N132").Value = -23958.4997
# Row 133
$ws.Range("H133").Value = 7298.75
$ws.Range("I133").Value = 3792.5
$ws.Range("J133").Value = 8000
$ws.Range("K133").Value = 11377.5
$ws.Range("L133").Value = 24000
$ws.Range("M133").Value = -6317.5
$ws.Range("N133").Value = -34120
# Row 135
$ws.Range("H135").Value = 1193.5333
$ws.Range("I135").Value = 260.42856
$ws.Range("J135").Value = 2010
$ws.Range("K135").Value = 2343.85704
$ws.Range("L135").Value = 18090
$ws.Range("M135").Value = 191.1429600000001
$ws.Range("N135").Value = -23160

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 1529.6471
$ws.Range("I102").Value = 1282.6364
$ws.Range("J102").Value = 1982.5
$ws.Range("K102").Value = 1282.6364
$ws.Range("L102").Value = 1982.5
$ws.Range("M102").Value = 339.3635999999999
$ws.Range("N102").Value = -5226.5
# Row 122
$ws.Range("H122").Value = 1615.7273
$ws.Range("I122").Value = 1501.375
$ws.Range("J122").Value = 1920.6666
$ws.Range("K122").Value = 4504.125
$ws.Range("L122").Value = 5761.9998
$ws.Range("M122").Value = -2054.125
$ws.Range("N122").Value = -10661.9998

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Range("H122").Value = 2602.3242
$ws.Range("I122").Value = 2364.6086
$ws.Range("J122").Value = 2992.8572
$ws.Range("K122").Value = 7093.825800000001
$ws.Range("L122").Value = 8978.571599999999
$ws.Range("M122").Value = -4643.825800000001
$ws.Range("N122").Value = -13878.5716

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 101391.914
$ws.Range("I122").Value = 1185.4286
$ws.Range("J122").Value = 241681
$ws.Range("K122").Value = 3556.2858
$ws.Range("L122").Value = 725043
$ws.Range("M122").Value = -1106.2858
$ws.Range("N122").Value = -729943
